# Update cfb_weather.xlsx with Timestamp 2025-09-14T10:01:35.816433
$wb = $excel.ActiveWorkbook

$wsFBS   = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- Refresh the run Timestamp column (FBS!AK2:AK45) ---
$wsFBS.Range("AK2:AK45").Value = "2025-09-14T10:01:35.816433"

# --- FBS sheet: wind_dir_fg (column Q) updates ---
$wsFBS.Range("Q10").Value = "NE"
$wsFBS.Range("Q13").Value = "W"
$wsFBS.Range("Q14").Value = "SSE"
$wsFBS.Range("Q17").Value = "W"
$wsFBS.Range("Q22").Value = "SSW"
$wsFBS.Range("Q24").Value = "SE"
$wsFBS.Range("Q29").Value = "ESE"
$wsFBS.Range("Q37").Value = "E"
$wsFBS.Range("Q39").Value = "SSW"
$wsFBS.Range("Q42").Value = "W"
$wsFBS.Range("Q43").Value = "W"
$wsFBS.Range("Q45").Value = "W"

# --- Other sheet: forecast temp/wind updates for row 2 ---
$wsOther.Range("Q2").Value = 79.7
$wsOther.Range("R2").Value = 8.699999999999999
$wsOther.Range("W2").Value = -5.3

# --- Other sheet: wind_dir_fg (column S) updates ---
$wsOther.Range("S14").Value = "SSE"
$wsOther.Range("S17").Value = "W"
$wsOther.Range("S23").Value = "N"
$wsOther.Range("S27").Value = "SSW"
$wsOther.Range("S31").Value = "ESE"
$wsOther.Range("S44").Value = "SSW"
$wsOther.Range("S47").Value = "W"
$wsOther.Range("S50").Value = "W"
